# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet gains three new columns describing how the
# record was produced by the scraper/importer:
#   - "category"    inserted right after "property_category"
#   - "source_file" appended after "legislator_id"
#   - "index"       appended after "source_file"
#
# New column layout (stock sheet):
#   B name | C owner | D quantity | E face_value | F currency | G total |
#   H property_category | I category | J date | K legislator_name |
#   L legislator_id | M source_file | N index

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the old "date" column (I) -- this shifts the
# existing date / legislator_name / legislator_id columns one place to the
# right (I->J, J->K, K->L) while keeping all of their values intact.
$ws.Columns("I:I").Insert()

# New "category" column header + values (order of assignment matters so
# that new shared strings are appended to the table in the expected order).
$ws.Range("I1").Value = "category"

# Two brand new columns appended at the end of the table.
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data rows: every stock record in this workbook came from the same
# normal/"normal" category, the same source file, and keeps its original
# row index.
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"

$ws.Range("M2").Value = "tmp77961"
$ws.Range("M3").Value = "tmp77961"
$ws.Range("M4").Value = "tmp77961"

$ws.Range("N2").Value = 75
$ws.Range("N3").Value = 76
$ws.Range("N4").Value = 77
